$p = $ppt.ActivePresentation

# --- Update the "datetimeFigureOut" date placeholder text on the slide master ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $shp.TextFrame.TextRange.Text = "02/08/2023"
    }
}

# --- Update the same date placeholder text on every slide layout ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "02/08/2023"
        }
    }
}

# --- Resize and re-label the "Target audience(s) in campaigns" textbox on slide 1 ---
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Target audience(s) in campaigns") {
        $shp.Width = 176.2373276346454
        $shp.TextFrame.TextRange.Text = "Target audience(s) in campaigns and journeys"
    }
}
